$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/15/2025  Through  9/21/2025"

# --- C14: was numeric 1, becomes text "0" (matches style/format of D14) ---
$ws.Range("D14").Copy($ws.Range("C14"))

# --- C28: was text "0", becomes numeric 2 (matches style/format of a numeric cell like C16) ---
$ws.Range("C16").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2

# --- Bulk numeric updates for rows 16-30 ---
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 76
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 26.666666666666
$ws.Range("L16").Value = -11.627906976744
$ws.Range("M16").Value = -53.658536585365
$ws.Range("N16").Value = -83.655913978494
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 46.153846153846
$ws.Range("I17").Value = 155
$ws.Range("J17").Value = 157
$ws.Range("K17").Value = -1.273885350318
$ws.Range("L17").Value = -11.428571428571
$ws.Range("M17").Value = 23.015873015873
$ws.Range("N17").Value = -19.270833333333
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 150
$ws.Range("I18").Value = 37
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -2.631578947368
$ws.Range("M18").Value = -79.891304347826
$ws.Range("N18").Value = -91.991341991342
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 75
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 34.782608695652
$ws.Range("I19").Value = 212
$ws.Range("J19").Value = 218
$ws.Range("K19").Value = -2.752293577981
$ws.Range("L19").Value = 36.774193548387
$ws.Range("M19").Value = 26.946107784431
$ws.Range("N19").Value = -26.896551724137
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 350
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 69.230769230769
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 105
$ws.Range("K20").Value = -4.761904761904
$ws.Range("L20").Value = -9.909909909909
$ws.Range("M20").Value = 23.456790123456
$ws.Range("N20").Value = -93.472584856396
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 41.176470588235
$ws.Range("F21").Value = 93
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = 47.619047619047
$ws.Range("I21").Value = 600
$ws.Range("J21").Value = 590
$ws.Range("K21").Value = 1.694915254237
$ws.Range("L21").Value = 3.270223752151
$ws.Range("M21").Value = -19.137466307277
$ws.Range("N21").Value = -79.770734996628
$ws.Range("F22").Value = 1
$ws.Range("L22").Value = -62.5
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -22.222222222222
$ws.Range("I23").Value = 71
$ws.Range("J23").Value = 74
$ws.Range("K23").Value = -4.054054054054
$ws.Range("L23").Value = -20.224719101123
$ws.Range("M23").Value = 69.047619047619
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 32
$ws.Range("H24").Value = 140.625
$ws.Range("I24").Value = 401
$ws.Range("J24").Value = 381
$ws.Range("K24").Value = 5.249343832021
$ws.Range("L24").Value = -4.066985645933
$ws.Range("M24").Value = -0.742574257425
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 17
$ws.Range("H25").Value = 183.333333333333
$ws.Range("I25").Value = 97
$ws.Range("J25").Value = 82
$ws.Range("K25").Value = 18.292682926829
$ws.Range("L25").Value = -11.009174311926
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = -23.076923076923
$ws.Range("I26").Value = 238
$ws.Range("J26").Value = 253
$ws.Range("K26").Value = -5.92885375494
$ws.Range("L26").Value = -11.524163568773
$ws.Range("M26").Value = -32.193732193732
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 21
$ws.Range("K28").Value = -12.5
$ws.Range("L28").Value = -4.545454545454
$ws.Range("M29").Value = -65.217391304347
$ws.Range("M30").Value = -65
